$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (column F) values to reflect repulled data
$ws.Range("F4").Value = -7
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -11
